$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its original "text" representation (values like
# "26.306.06" or "0.9980" must not be reinterpreted/normalized as numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.306.06'
$ws.Range("E2").Value = '  +4.73%  '
$ws.Range("D3").Value = '1.714.57'
$ws.Range("E3").Value = '  +3.96%  '
$ws.Range("D4").Value = '0.9980'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '240.97'
$ws.Range("E5").Value = '  +3.12%  '
$ws.Range("D6").Value = '0.9985'
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").Value = '0.4719'
$ws.Range("E7").Value = '  -0.79%  '
$ws.Range("D8").Value = '0.2646'
$ws.Range("E8").Value = '  +3.22%  '
$ws.Range("D9").Value = '0.06243'
$ws.Range("E9").Value = '  +2.30%  '
$ws.Range("D10").Value = '1.697.69'
$ws.Range("E10").Value = '  +2.91%  '
$ws.Range("D11").Value = '0.07079'
$ws.Range("E11").Value = '  +0.64%  '
$ws.Range("D12").Value = '15.28'
$ws.Range("E12").Value = '  +6.46%  '
$ws.Range("D13").Value = '0.5906'
$ws.Range("E13").Value = '  +2.74%  '
$ws.Range("D14").Value = '4.417'
$ws.Range("E14").Value = '  +2.66%  '
$ws.Range("D15").Value = '76.15'
$ws.Range("E15").Value = '  +3.52%  '
$ws.Range("D16").Value = '0.9985'
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("D17").Value = '0.9990'
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").Value = '26.315.11'
$ws.Range("E18").Value = '  +4.82%  '
$ws.Range("D19").Value = '0.000006810'
$ws.Range("E19").Value = '  +2.77%  '
$ws.Range("D20").Value = '11.62'
$ws.Range("E20").Value = '  +2.72%  '
$ws.Range("D21").Value = '1.917.12'
$ws.Range("E21").Value = '  +3.58%  '
$ws.Range("D22").Value = '4.575'
$ws.Range("E22").Value = '  +6.00%  '
$ws.Range("D23").Value = '8.852'
$ws.Range("E23").Value = '  +4.86%  '
$ws.Range("D24").Value = '5.343'
$ws.Range("E24").Value = '  +1.82%  '
$ws.Range("D25").Value = '135.47'
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("D26").Value = '15.18'
$ws.Range("E26").Value = '  +1.63%  '
$ws.Range("D27").Value = '1.403'
$ws.Range("E27").Value = '  +2.15%  '
$ws.Range("D28").Value = '1.762'
$ws.Range("E28").Value = '  +7.65%  '
$ws.Range("D29").Value = '106.48'
$ws.Range("E29").Value = '  +2.98%  '
$ws.Range("D30").Value = '4.025'
$ws.Range("E30").Value = '  +3.34%  '
$ws.Range("D31").Value = '3.700'
$ws.Range("E31").Value = '  +4.86%  '
$ws.Range("D32").Value = '0.07771'
$ws.Range("E32").Value = '  +2.35%  '
$ws.Range("D33").Value = '0.04422'
$ws.Range("E33").Value = '  +3.96%  '
$ws.Range("D34").Value = '2.610'
$ws.Range("E34").Value = '  +1.35%  '
$ws.Range("D35").Value = '0.6226'
$ws.Range("E35").Value = '  +4.99%  '
$ws.Range("D36").Value = '0.9730'
$ws.Range("E36").Value = '  +4.05%  '
$ws.Range("D37").Value = '0.9171'
$ws.Range("E37").Value = '  +7.65%  '
$ws.Range("D38").Value = '111.36'
$ws.Range("E38").Value = '  +13.15%  '
$ws.Range("D39").Value = '2.402'
$ws.Range("E39").Value = '  -7.11%  '
$ws.Range("D40").Value = '1.916'
$ws.Range("E40").Value = '  +7.50%  '
$ws.Range("D41").Value = '0.9996'
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = '0.01470'
$ws.Range("E42").Value = '  -0.64%  '
$ws.Range("D43").Value = '0.3815'
$ws.Range("E43").Value = '  +3.68%  '
$ws.Range("D44").Value = '5.137'
$ws.Range("E44").Value = '  +11.00%  '
$ws.Range("D45").Value = '0.1142'
$ws.Range("E45").Value = '  +4.03%  '
$ws.Range("D46").Value = '6.256'
$ws.Range("E46").Value = '  +2.98%  '
$ws.Range("D47").Value = '0.05301'
$ws.Range("E47").Value = '  +1.58%  '
$ws.Range("D48").Value = '30.74'
$ws.Range("E48").Value = '  +6.03%  '
$ws.Range("D49").Value = '7.663'
$ws.Range("E49").Value = '  +7.05%  '
$ws.Range("D50").Value = '1.224'
$ws.Range("E50").Value = '  +2.40%  '
$ws.Range("D51").Value = '0.3384'
$ws.Range("E51").Value = '  +3.20%  '
